$d = $word.ActiveDocument

# --- Paragraph 1 (Heading1): center alignment + bold/black/size-28 run formatting ---
$p1 = $d.Paragraphs.Item(1)
$p1.Alignment = 1

$r1 = $p1.Range
$r1run = $d.Range($r1.Start, $r1.End - 1)
$r1run.Font.Bold = 1
$r1run.Font.Color = 0
$r1run.Font.Size = 14

# --- Paragraph 2: move leading <w:br/> to the very start, add 4-space indent and
#     xml:space="preserve" to every line, and insert a new break before the last line ---
$p2 = $d.Paragraphs.Item(2)
$r2 = $p2.Range

$vt = [char]11

$line1 = "Характеристики ИЗАВ и показатели выбросов определяются для всех основных режимов работы технологического оборудования (установок) и стадий технологических процессов."
$line2 = "В ходе инвентаризации выбросов при определении качественных и количественных показателей выбросов выявлены, учтены и проанализированы изменения показателей выбросов во времени, обусловленные неодновременной, неравномерной работой оборудования, изменениями режимов работы оборудования и стадийностью процессов, в ходе которых образуются и выделяются загрязняющие вещества (табл. 3.1, 3.2)."
$line3 = "Суммарные выбросы ЗВ в атмосферный воздух с учетом их очистки и утилизации (в целом по объекту ОНВ) учтены в таблице 3.7."
$line4 = "Выбросы от передвижных ИЗАВ учтены в таблице 3.8."

$newText = $vt + "    " + $line1 + $vt + "    " + $line2 + $vt + "    " + $line3 + $vt + "    " + $line4

$r2.Text = $newText

Write-Host "Edit applied"
